$d = $word.ActiveDocument

$replacements = @(
    @{old="601÷7=85, 6"; new="406÷5=81, 1"},
    @{old="703÷6=117, 1"; new="810÷3=270, 0"},
    @{old="931÷2=465, 1"; new="549÷8=68, 5"},
    @{old="842÷3=280, 2"; new="212÷6=35, 2"},
    @{old="235÷5=47, 0"; new="917÷9=101, 8"},
    @{old="163÷9=18, 1"; new="273÷4=68, 1"},
    @{old="174÷2=87, 0"; new="672÷2=336, 0"},
    @{old="137÷3=45, 2"; new="707÷8=88, 3"},
    @{old="994÷4=248, 2"; new="363÷8=45, 3"},
    @{old="141÷7=20, 1"; new="311÷7=44, 3"},
    @{old="738÷8=92, 2"; new="625÷4=156, 1"},
    @{old="956÷4=239, 0"; new="991÷7=141, 4"},
    @{old="406÷7=58, 0"; new="295÷8=36, 7"},
    @{old="970÷5=194, 0"; new="858÷4=214, 2"},
    @{old="973÷9=108, 1"; new="633÷5=126, 3"},
    @{old="784÷2=392, 0"; new="803÷5=160, 3"},
    @{old="328÷4=82, 0"; new="144÷4=36, 0"},
    @{old="417÷3=139, 0"; new="479÷8=59, 7"},
    @{old="822÷5=164, 2"; new="480÷2=240, 0"},
    @{old="669÷9=74, 3"; new="607÷9=67, 4"},
    @{old="763÷2=381, 1"; new="321÷4=80, 1"},
    @{old="699÷7=99, 6"; new="168÷5=33, 3"},
    @{old="473÷6=78, 5"; new="202÷8=25, 2"},
    @{old="641÷7=91, 4"; new="108÷5=21, 3"},
    @{old="974÷3=324, 2"; new="966÷4=241, 2"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
